# Add two new columns (I: "I0", J: "IF") to the sheet, mirroring the
# header style used by the existing "IP" header in H1, and fill in the
# per-row numeric values for rows 2-24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers - copy the style from H1 (bold/centered/bordered) onto I1:J1
$ws.Range("H1").Copy($ws.Range("I1:J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (col 9) and J (col 10), rows 2-24
$values = @{
    2  = @(8, 8)
    3  = @(6, 6)
    4  = @(1, 2)
    5  = @(3, 5)
    6  = @(1, 4)
    7  = @(7, 7)
    8  = @(5, 7)
    9  = @(4, 5)
    10 = @(7, 7)
    11 = @(8, 8)
    12 = @(2, 4)
    13 = @(6, 6)
    14 = @(3, 5)
    15 = @(5, 7)
    16 = @(4, 5)
    17 = @(4, 4)
    18 = @(4, 4)
    19 = @(7, 7)
    20 = @(7, 8)
    21 = @(8, 9)
    22 = @(4, 5)
    23 = @(5, 5)
    24 = @(7, 7)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
